$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.734.11"
$ws.Range("E2").Value = "  -3.53%  "

$ws.Range("D3").Value = "2.400.24"
$ws.Range("E3").Value = "  -4.03%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'505.64"
$ws.Range("E5").Value = "  -5.19%  "

$ws.Range("D6").Value = "'130.09"
$ws.Range("E6").Value = "  -2.71%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("E8").Value = "  -2.87%  "

$ws.Range("D9").Value = "2.398.86"
$ws.Range("E9").Value = "  -4.22%  "

$ws.Range("D10").Value = "'0.0965"
$ws.Range("E10").Value = "  -2.75%  "

$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "  -1.39%  "

$ws.Range("D12").Value = "'0.322"
$ws.Range("E12").Value = "  -2.13%  "

$ws.Range("D13").Value = "'4.67"
$ws.Range("E13").Value = "  -10.14%  "

$ws.Range("D14").Value = "2.798.95"
$ws.Range("E14").Value = "  -4.68%  "

$ws.Range("D15").Value = "56.508.16"
$ws.Range("E15").Value = "  -3.66%  "

$ws.Range("D16").Value = "'21.64"
$ws.Range("E16").Value = "  -2.84%  "

$ws.Range("E17").Value = "  -3.00%  "

$ws.Range("D18").Value = "2.382.21"
$ws.Range("E18").Value = "  -4.57%  "

$ws.Range("D19").Value = "'10.24"
$ws.Range("E19").Value = "  -3.25%  "

$ws.Range("E20").Value = "  -2.54%  "

$ws.Range("E21").Value = "  -4.39%  "

$ws.Range("E22").Value = "  +0.46%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("D24").Value = "'65.71"
$ws.Range("E24").Value = "  -0.44%  "

$ws.Range("D25").Value = "'0.996"
$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("D26").Value = "2.490.69"
$ws.Range("E26").Value = "  -4.54%  "

$ws.Range("D27").Value = "'0.376"
$ws.Range("E27").Value = "  -7.71%  "

$ws.Range("E28").Value = "  -4.95%  "

$ws.Range("E29").Value = "  -2.41%  "

$ws.Range("D30").Value = "'174.93"
$ws.Range("E30").Value = "  +1.04%  "

$ws.Range("E31").Value = "  -3.51%  "

$ws.Range("E32").Value = "  -5.48%  "

$ws.Range("D33").Value = "'6.16"
$ws.Range("E33").Value = "  -2.00%  "

$ws.Range("E34").Value = "  -5.76%  "

$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("D37").Value = "'17.83"
$ws.Range("E37").Value = "  -1.35%  "

$ws.Range("E38").Value = "  -0.62%  "

$ws.Range("E39").Value = "  -4.51%  "

$ws.Range("E40").Value = "  -1.41%  "

$ws.Range("E41").Value = "  -4.68%  "

$ws.Range("D42").Value = "'0.789"
$ws.Range("E42").Value = "  -5.51%  "

$ws.Range("D43").Value = "'131.58"
$ws.Range("E43").Value = "  +0.33%  "

$ws.Range("E44").Value = "  -2.70%  "

$ws.Range("D45").Value = "'4.89"
$ws.Range("E45").Value = "  -2.29%  "

$ws.Range("D46").Value = "'255.55"
$ws.Range("E46").Value = "  -7.04%  "

$ws.Range("D47").Value = "'0.571"
$ws.Range("E47").Value = "  -3.48%  "

$ws.Range("E48").Value = "  -3.44%  "

$ws.Range("E49").Value = "  -3.83%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'16.89"
$ws.Range("E50").Value = "  -3.77%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0208"
$ws.Range("E51").Value = "  -4.42%  "
